# Insert two new data rows (313 and 314) into the "Arándano (blue)" sheet,
# pushing the existing rows 313-377 down to 315-379.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows before the current row 313 (shifts 313..377 -> 315..379)
$ws.Range("A313:T314").EntireRow.Insert()

# --- New row 313 ---
$row313 = @(9, "Vega Central Mapocho de Santiago", "Metropolitana", 45258, 13, "Fruta", 100101, "Berries", 100101001, "Arándano (blue)", "Sin especificar", "Especial", 500, 4600, 5000, 4800, "`$/bandeja 2 kilos", "Provincia de Curicó", 2400, 2)
for ($i = 0; $i -lt $row313.Length; $i++) {
    $ws.Cells.Item(313, $i + 1).Value = $row313[$i]
}

# --- New row 314 ---
$row314 = @(9, "Vega Central Mapocho de Santiago", "Metropolitana", 45258, 13, "Fruta", 100101, "Berries", 100101001, "Arándano (blue)", "Sin especificar", "Primera", 600, 4000, 4000, 4000, "`$/bandeja 2 kilos", "Provincia de Curicó", 2000, 2)
for ($i = 0; $i -lt $row314.Length; $i++) {
    $ws.Cells.Item(314, $i + 1).Value = $row314[$i]
}
